# Weekly update of fruit/vegetable price data (Pepino dulce - Agricola del Norte S.A. de Arica)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44363
$ws.Range("H2").Value = 'Cultivar IV Región'
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14500
$ws.Range("N2").Value = '$/bandeja 18 kilos'
$ws.Range("O2").Value = 'Provincia de Limarí'
$ws.Range("P2").Value = 806
$ws.Range("Q2").Value = 18

# Row 3
$ws.Range("D3").Value = 44391
$ws.Range("I3").Value = 'Segunda'
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 15500
$ws.Range("P3").Value = 861

# Row 4
$ws.Range("D4").Value = 44435
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17500
$ws.Range("P4").Value = 972

# Row 5
$ws.Range("D5").Value = 44435
$ws.Range("I5").Value = 'Tercera'
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("P5").Value = 806

# Row 6
$ws.Range("D6").Value = 44412
$ws.Range("H6").Value = 'Cultivar IV Región'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17500
$ws.Range("N6").Value = '$/bandeja 18 kilos'
$ws.Range("O6").Value = 'Provincia de Limarí'
$ws.Range("P6").Value = 972
$ws.Range("Q6").Value = 18

# Row 7
$ws.Range("D7").Value = 44377
$ws.Range("J7").Value = 100
$ws.Range("M7").Value = 17600
$ws.Range("P7").Value = 978

# Row 8
$ws.Range("D8").Value = 44405
$ws.Range("I8").Value = 'Segunda'
$ws.Range("J8").Value = 140
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17500
$ws.Range("P8").Value = 972

# Row 9
$ws.Range("D9").Value = 44398
$ws.Range("I9").Value = 'Primera'
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 17500
$ws.Range("P9").Value = 972

# Row 10
$ws.Range("D10").Value = 44398
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15500
$ws.Range("P10").Value = 861

# Row 11
$ws.Range("D11").Value = 44454
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 19000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 19500
$ws.Range("P11").Value = 1083

# Row 12
$ws.Range("D12").Value = 44221
$ws.Range("H12").Value = 'Cultivar XV región'
$ws.Range("J12").Value = 140
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = 5500
$ws.Range("N12").Value = '$/caja 10 kilos'
$ws.Range("O12").Value = 'Región de Arica y Parinacota'
$ws.Range("P12").Value = 550
$ws.Range("Q12").Value = 10

# Row 13
$ws.Range("D13").Value = 44211
$ws.Range("H13").Value = 'Cultivar XV región'
$ws.Range("J13").Value = 140
$ws.Range("K13").Value = 4500
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 4750
$ws.Range("N13").Value = '$/caja 10 kilos'
$ws.Range("O13").Value = 'Región de Arica y Parinacota'
$ws.Range("P13").Value = 475
$ws.Range("Q13").Value = 10

# Row 14
$ws.Range("I14").Value = 'Segunda'
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 17000
$ws.Range("L14").Value = 18000
$ws.Range("M14").Value = 17500
$ws.Range("P14").Value = 972

# Row 15
$ws.Range("D15").Value = 44433
$ws.Range("I15").Value = 'Tercera'
$ws.Range("J15").Value = 120
